$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Users_OnGoing" (sheet2.xml): fill in the Partner Farmer profile
# credentials (row 12 = username w/ mailto hyperlink, row 13 = password).
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users_OnGoing")
$wsUsers.Range("B12").Value = "soi.testing.crew@gmail.com.farmer"
$wsUsers.Hyperlinks.Add($wsUsers.Range("B12"), "mailto:soi.testing.crew@gmail.com.farmer") | Out-Null
$wsUsers.Range("B12").Style = "Hyperlink"
$wsUsers.Range("B13").Value = "Testing-2020"
$wsUsers.Range("B13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Environments_OnGoing" (sheet1.xml): append two new variable rows.
# ---------------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("Environments_OnGoing")
$wsEnv.Range("B12").Value = "/s/opportunity/related"
$wsEnv.Range("A12").Value = "PartnersCommunityOpportunityRelated"
$wsEnv.Range("A13").Value = "PartnersCommunityOpportunitesSufix"
$wsEnv.Range("B13").Value = "/Opportunities"

# ---------------------------------------------------------------------------
# Sheet "Companies" (sheet3.xml): only the saved view (scroll/selection)
# changed. The sheet needs to be activated momentarily so its window state
# (scroll position + selection) can be updated, then we restore the
# originally active sheet so "tabSelected" stays where it was.
# ---------------------------------------------------------------------------
$wsCompanies = $wb.Worksheets.Item("Companies")
$wsCompanies.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$wsCompanies.Range("B38").Select() | Out-Null

# Restore the originally active sheet/selection.
$wsEnv.Activate()
$wsEnv.Range("B17").Select() | Out-Null
